# Bump the "Förändrad" (Changed) date column (C) by one day.
#
# Every data row (2..last) currently stores the same serial date value
# 45171 (2023-09-02) in column C; the update moves it to 45172
# (2023-09-03) for every one of those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlUp = -4162
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End($xlUp).Row
if ($lastRow -lt 2) { $lastRow = 397 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45171) {
        $cell.Value2 = 45172
    }
}
